$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above row 179, pushing the existing
# rows 179-208 down to 181-210 (dimension grows from A1:T208 to A1:T210).
$ws.Rows.Item(179).Insert()
$ws.Rows.Item(179).Insert()

# Populate the two new rows (179-180) with this week's price data.
# Row 179: Cereza, Royal Dawn, Primera - Provincia de Curicó
$ws.Cells.Item(179, 1).Value = 11
$ws.Cells.Item(179, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(179, 3).Value = "Bíobío"
$ws.Cells.Item(179, 4).Value = 45265
$ws.Cells.Item(179, 5).Value = 8
$ws.Cells.Item(179, 6).Value = "Fruta"
$ws.Cells.Item(179, 7).Value = 100103
$ws.Cells.Item(179, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(179, 9).Value = 100103001
$ws.Cells.Item(179, 10).Value = "Cereza"
$ws.Cells.Item(179, 11).Value = "Royal Dawn"
$ws.Cells.Item(179, 12).Value = "Primera"
$ws.Cells.Item(179, 13).Value = 150
$ws.Cells.Item(179, 14).Value = 12000
$ws.Cells.Item(179, 15).Value = 12000
$ws.Cells.Item(179, 16).Value = 12000
$ws.Cells.Item(179, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(179, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(179, 19).Value = 1200
$ws.Cells.Item(179, 20).Value = 10

# Row 180: Cereza, Santina, Primera - Provincia de Curicó
$ws.Cells.Item(180, 1).Value = 11
$ws.Cells.Item(180, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(180, 3).Value = "Bíobío"
$ws.Cells.Item(180, 4).Value = 45265
$ws.Cells.Item(180, 5).Value = 8
$ws.Cells.Item(180, 6).Value = "Fruta"
$ws.Cells.Item(180, 7).Value = 100103
$ws.Cells.Item(180, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(180, 9).Value = 100103001
$ws.Cells.Item(180, 10).Value = "Cereza"
$ws.Cells.Item(180, 11).Value = "Santina"
$ws.Cells.Item(180, 12).Value = "Primera"
$ws.Cells.Item(180, 13).Value = 120
$ws.Cells.Item(180, 14).Value = 13000
$ws.Cells.Item(180, 15).Value = 13000
$ws.Cells.Item(180, 16).Value = 13000
$ws.Cells.Item(180, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(180, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(180, 19).Value = 1300
$ws.Cells.Item(180, 20).Value = 10
